$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Update publish-date labels ---
# Row 9 holds the "publish date" for each reporting column; the dates used
# for the newest quarter (column M, "فصل دوم منتهی به 1401/12") move
# forward from 1402-02-23 to 1402-03-09, with incremented revision counts.
$ws.Range("I9").Value = "1402-03-09 (6)"
$ws.Range("K9").Value = "1402-03-09 (6)"
$ws.Range("M9").Value = "1402-03-09 (3)"

# --- Update the restated income-statement figures for column M (the
#     "فصل دوم منتهی به 1401/12" quarter) ---
$ws.Range("M12").Value = -16914
$ws.Range("M13").Value = 8439
$ws.Range("M14").Value = -695
$ws.Range("M17").Value = 8769
$ws.Range("M20").Value = 8267
$ws.Range("M21").Value = -1788
$ws.Range("M22").Value = 6479
$ws.Range("M24").Value = 6479
